# "Generate Report for Handoff"
# Adds a new handed-off file ("2281ceb4-...") as row 3 on the Overview,
# zh-cn and de-de sheets/tables, mirroring the existing 6e4d8e53-... row.

$wb = $excel.ActiveWorkbook

# ---- build the long repeated-"o" filler filenames -------------------------
$guidNew   = "2281ceb4-b1b9-428b-958d-8229bb1a6188"
$commitNew = "af4e871ade7abd3eeb1fdf9421f80be5e817d5ea"

$fillerMd  = "".PadRight(149, 'o')
$fillerXlf = "".PadRight(40, 'o')

$mdName      = $guidNew + $fillerMd + ".md"
$mdDisplay   = "e2e\" + $mdName
$zhXlfName   = $guidNew + $fillerXlf + "." + $commitNew + ".zh-cn.xlf"
$deXlfName   = $guidNew + $fillerXlf + "." + $commitNew + ".de-de.xlf"

$readyForHandoff = "Ready for handoff"
$hoDate          = "2016-08-29 06:28:32"
$handoffDate     = "2016-08-29 06:28:27"
$zeroDate        = "0001-01-01 00:00:00"

# GitHub blob URLs follow the same pattern as the existing hyperlink target.
$urlBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bdca7fbfd554df70405dce753df4d0e1bddf4de0/e2e/"
$mdUrl   = $urlBase + $mdName

# =========================== Overview sheet =================================
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item(1)
$loOv.ListRows.Add() | Out-Null

$wsOv.Range("A3").Value = $mdName
$wsOv.Range("B3").Value = $mdDisplay
$wsOv.Range("C3").Value = ".md"
$wsOv.Range("D3").Value = ""
$wsOv.Range("E3").Value = $readyForHandoff
$wsOv.Range("F3").Value = $readyForHandoff
$wsOv.Range("G3").Value = $hoDate

$wsOv.Hyperlinks.Add($wsOv.Range("B3"), $mdUrl, "", "", $mdDisplay) | Out-Null

$wsOv.Columns.Item(5).ColumnWidth = 16.38265482584637
$wsOv.Columns.Item(6).ColumnWidth = 16.38265482584637

# ============================ zh-cn sheet ====================================
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A3").Value = $mdName
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $readyForHandoff
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'False"
$wsZh.Range("G3").Value = $zhXlfName
$wsZh.Range("H3").Value = $handoffDate
$wsZh.Range("K3").Value = $zeroDate
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("O3").Value = "'False"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrl, "", "", $mdName) | Out-Null

$wsZh.Columns.Item(3).ColumnWidth = 16.38265482584637

# ============================ de-de sheet ====================================
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A3").Value = $mdName
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $readyForHandoff
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'False"
$wsDe.Range("G3").Value = $deXlfName
$wsDe.Range("H3").Value = $hoDate
$wsDe.Range("K3").Value = $zeroDate
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("O3").Value = "'False"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrl, "", "", $mdName) | Out-Null

$wsDe.Columns.Item(3).ColumnWidth = 16.38265482584637

Write-Output "Report row added to Overview, zh-cn and de-de sheets."
